$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.051248330555805
$ws.Range("D2").Value = 1.057180132492496
$ws.Range("E2").Value = 1.048143192356557
$ws.Range("F2").Value = 1.064939561016133
$ws.Range("I2").Value = 1.049168331995405
$ws.Range("J2").Value = 1.05627720495258
$ws.Range("K2").Value = 1.059915808933184
$ws.Range("L2").Value = 1.050903860449022
$ws.Range("M2").Value = 1.067654145098621
$ws.Range("N2").Value = 1.057777239740302

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05265992944496
$ws.Range("D3").Value = 1.058068084271872
$ws.Range("E3").Value = 1.049365062568878
$ws.Range("F3").Value = 1.066091568972227
$ws.Range("I3").Value = 1.04960547009556
$ws.Range("J3").Value = 1.057336558526316
$ws.Range("K3").Value = 1.060617357570881
$ws.Range("L3").Value = 1.051936655679465
$ws.Range("M3").Value = 1.068620613010205
$ws.Range("N3").Value = 1.058838097717622

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053572276119186
$ws.Range("D4").Value = 1.058641961211886
$ws.Range("E4").Value = 1.050154849252077
$ws.Range("F4").Value = 1.0668365929837
$ws.Range("I4").Value = 1.049886531038929
$ws.Range("J4").Value = 1.058020538702663
$ws.Range("K4").Value = 1.061069983700112
$ws.Range("L4").Value = 1.052603537049207
$ws.Range("M4").Value = 1.069244997405348
$ws.Range("N4").Value = 1.059523049224273

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053955580140386
$ws.Range("D5").Value = 1.05888305618524
$ws.Range("E5").Value = 1.050486677166228
$ws.Range("F5").Value = 1.067149707368734
$ws.Range("I5").Value = 1.050004260476436
$ws.Range("J5").Value = 1.058307730436095
$ws.Range("K5").Value = 1.061259952790646
$ws.Range("L5").Value = 1.052883560905975
$ws.Range("M5").Value = 1.069507255212516
$ws.Range("N5").Value = 1.059810648802908

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.054019924272003
$ws.Range("D6").Value = 1.058923527587304
$ws.Range("E6").Value = 1.050542381011441
$ws.Range("F6").Value = 1.067202275231042
$ws.Range("I6").Value = 1.050024002679189
$ws.Range("J6").Value = 1.058355930569384
$ws.Range("K6").Value = 1.061291831019183
$ws.Range("L6").Value = 1.05293055871594
$ws.Range("M6").Value = 1.069551275835909
$ws.Range("N6").Value = 1.059858917385915

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.05357739880826
$ws.Range("D7").Value = 1.058645183372759
$ws.Range("E7").Value = 1.050159283927676
$ws.Range("F7").Value = 1.066840777198464
$ws.Range("I7").Value = 1.049888105827576
$ws.Range("J7").Value = 1.058024377557151
$ws.Range("K7").Value = 1.06107252331173
$ws.Range("L7").Value = 1.052607280045997
$ws.Range("M7").Value = 1.069248502621092
$ws.Range("N7").Value = 1.059526893530375

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.051725608116478
$ws.Range("D8").Value = 1.057480363237988
$ws.Range("E8").Value = 1.048556306439736
$ws.Range("F8").Value = 1.065328972321889
$ws.Range("I8").Value = 1.049316437795605
$ws.Range("J8").Value = 1.05663552973436
$ws.Range("K8").Value = 1.060153175511816
$ws.Range("L8").Value = 1.051253191754251
$ws.Range("M8").Value = 1.067980972854145
$ws.Range("N8").Value = 1.058136073384372

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.048454216399903
$ws.Range("D9").Value = 1.055422463925075
$ws.Range("E9").Value = 1.045725008060489
$ws.Range("F9").Value = 1.062661782932721
$ws.Range("I9").Value = 1.048295256380686
$ws.Range("J9").Value = 1.054176603435648
$ws.Range("K9").Value = 1.058522958199082
$ws.Range("L9").Value = 1.048856184588174
$ws.Range("M9").Value = 1.065739771066025
$ws.Range("N9").Value = 1.055673655128273

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.046267369000302
$ws.Range("D10").Value = 1.054046836301113
$ws.Range("E10").Value = 1.043832759588853
$ws.Range("F10").Value = 1.060881343056227
$ws.Range("I10").Value = 1.047605072551103
$ws.Range("J10").Value = 1.052529282754634
$ws.Range("K10").Value = 1.057429171031676
$ws.Range("L10").Value = 1.0472506041001
$ws.Range("M10").Value = 1.06424033708536
$ws.Range("N10").Value = 1.054023995062942

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045318959878337
$ws.Range("D11").Value = 1.053450273545306
$ws.Range("E11").Value = 1.043012220894478
$ws.Range("F11").Value = 1.060109806672348
$ws.Range("I11").Value = 1.047303964377753
$ws.Range("J11").Value = 1.05181401501096
$ws.Range("K11").Value = 1.056953869081788
$ws.Range("L11").Value = 1.046553522166213
$ws.Range("M11").Value = 1.063589773544291
$ws.Range("N11").Value = 1.05330771155703

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044966448052417
$ws.Range("D12").Value = 1.053228545116226
$ws.Range("E12").Value = 1.042707253738071
$ws.Range("F12").Value = 1.059823131109918
$ws.Range("I12").Value = 1.047191778762562
$ws.Range("J12").Value = 1.051548032498836
$ws.Range("K12").Value = 1.056777065354286
$ws.Range("L12").Value = 1.046294311728169
$ws.Range("M12").Value = 1.063347927038186
$ws.Range("N12").Value = 1.05304135131923

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045042073570904
$ws.Range("D13").Value = 1.053276112918155
$ws.Range("E13").Value = 1.042772678567761
$ws.Range("F13").Value = 1.059884628235126
$ws.Range("I13").Value = 1.047215858404311
$ws.Range("J13").Value = 1.0516051003122
$ws.Range("K13").Value = 1.056815001949199
$ws.Range("L13").Value = 1.046349926133289
$ws.Range("M13").Value = 1.063399812955864
$ws.Range("N13").Value = 1.053098500175434

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045289825880431
$ws.Range("D14").Value = 1.053431948241553
$ws.Range("E14").Value = 1.042987015966481
$ws.Range("F14").Value = 1.06008611188922
$ws.Range("I14").Value = 1.047294698040176
$ws.Range("J14").Value = 1.051792034965332
$ws.Range("K14").Value = 1.056939259654297
$ws.Range("L14").Value = 1.046532101553133
$ws.Range("M14").Value = 1.06356978649883
$ws.Range("N14").Value = 1.053285700297216

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045442443519363
$ws.Range("D15").Value = 1.053527945113952
$ws.Range("E15").Value = 1.043119051975758
$ws.Range("F15").Value = 1.060210240248317
$ws.Range("I15").Value = 1.047343228544923
$ws.Range("J15").Value = 1.051907171650888
$ws.Range("K15").Value = 1.057015785013089
$ws.Range("L15").Value = 1.046644308184601
$ws.Range("M15").Value = 1.063674486445951
$ws.Range("N15").Value = 1.053401000490071

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.046330279837527
$ws.Range("D16").Value = 1.05408640891098
$ws.Range("E16").Value = 1.043887190707687
$ws.Range("F16").Value = 1.060932534582929
$ws.Range("I16").Value = 1.047625008449566
$ws.Range("J16").Value = 1.052576710875545
$ws.Range("K16").Value = 1.057460679563407
$ws.Range("L16").Value = 1.047296827668074
$ws.Range("M16").Value = 1.064283485216333
$ws.Range("N16").Value = 1.054071490537225

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04688679258156
$ws.Range("D17").Value = 1.054436474137715
$ws.Range("E17").Value = 1.044368703785966
$ws.Range("F17").Value = 1.061385449316868
$ws.Range("I17").Value = 1.047801156671652
$ws.Range("J17").Value = 1.052996165204062
$ws.Range("K17").Value = 1.057739297460203
$ws.Range("L17").Value = 1.047705636113336
$ws.Range("M17").Value = 1.064665144151194
$ws.Range("N17").Value = 1.054491540539

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.047211253759273
$ws.Range("D18").Value = 1.054640573908488
$ws.Range("E18").Value = 1.044649448444036
$ws.Range("F18").Value = 1.061649569717047
$ws.Range("I18").Value = 1.047903683604793
$ws.Range("J18").Value = 1.053240636348415
$ws.Range("K18").Value = 1.057901648154251
$ws.Range("L18").Value = 1.047943908488243
$ws.Range("M18").Value = 1.06488763425472
$ws.Range("N18").Value = 1.054736358860419

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.047321862531153
$ws.Range("D19").Value = 1.054710151905314
$ws.Range("E19").Value = 1.044745155966015
$ws.Range("F19").Value = 1.061739618362586
$ws.Range("I19").Value = 1.047938605833896
$ws.Range("J19").Value = 1.053323962695243
$ws.Range("K19").Value = 1.057956978066192
$ws.Range("L19").Value = 1.048025123039039
$ws.Range("M19").Value = 1.064963476505615
$ws.Range("N19").Value = 1.054819803540218

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.046827098900281
$ws.Range("D20").Value = 1.054398924503357
$ws.Range("E20").Value = 1.044317053793304
$ws.Range("F20").Value = 1.061336861782815
$ws.Range("I20").Value = 1.047782280129951
$ws.Range("J20").Value = 1.052951181356139
$ws.Range("K20").Value = 1.057709421212325
$ws.Range("L20").Value = 1.047661793349831
$ws.Range("M20").Value = 1.064624208721488
$ws.Range("N20").Value = 1.054446492808853

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045216875415316
$ws.Range("D21").Value = 1.053386062462879
$ws.Range("E21").Value = 1.042923904023373
$ws.Range("F21").Value = 1.060026782557923
$ws.Range("I21").Value = 1.04727149115163
$ws.Range("J21").Value = 1.051736995694789
$ws.Range("K21").Value = 1.056902675926479
$ws.Range("L21").Value = 1.046478463285578
$ws.Range("M21").Value = 1.063519739024488
$ws.Range("N21").Value = 1.053230582864595

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044203127376864
$ws.Range("D22").Value = 1.052748433705983
$ws.Range("E22").Value = 1.04204691685029
$ws.Range("F22").Value = 1.059202546348942
$ws.Range("I22").Value = 1.046948366055733
$ws.Range("J22").Value = 1.050971849620571
$ws.Range("K22").Value = 1.056393964217062
$ws.Range("L22").Value = 1.045732816134757
$ws.Range("M22").Value = 1.062824166568752
$ws.Range("N22").Value = 1.052464350195188

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.044740663311437
$ws.Range("D23").Value = 1.053086529553793
$ws.Range("E23").Value = 1.042511926323257
$ws.Range("F23").Value = 1.059639541615234
$ws.Range("I23").Value = 1.047119848353832
$ws.Range("J23").Value = 1.051377634479259
$ws.Range("K23").Value = 1.056663782871554
$ws.Range("L23").Value = 1.046128254806466
$ws.Range("M23").Value = 1.063193012601728
$ws.Range("N23").Value = 1.05287071131492

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.046854072348476
$ws.Range("D24").Value = 1.054415891837106
$ws.Range("E24").Value = 1.044340392553661
$ws.Range("F24").Value = 1.061358816574893
$ws.Range("I24").Value = 1.047790810298373
$ws.Range("J24").Value = 1.05297150820647
$ws.Range("K24").Value = 1.05772292150486
$ws.Range("L24").Value = 1.047681604559876
$ws.Range("M24").Value = 1.064642706067026
$ws.Range("N24").Value = 1.054466848525644

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.049300967971932
$ws.Range("D25").Value = 1.055955123603149
$ws.Range("E25").Value = 1.046457779634014
$ws.Range("F25").Value = 1.063351710960077
$ws.Range("I25").Value = 1.048560905164986
$ws.Range("J25").Value = 1.054813693780952
$ws.Range("K25").Value = 1.058945629713692
$ws.Range("L25").Value = 1.049477186980676
$ws.Range("M25").Value = 1.066320098127007
$ws.Range("N25").Value = 1.056311650214942
